$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.371.86"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.879.07"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7207"
$ws.Range("E5").Value = "  +1.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.05"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08034"
$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3137"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  -0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08151"
$ws.Range("E11").Value = "  -3.23%  "

$ws.Range("D12").Value = "1.882.80"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.51"
$ws.Range("E13").Value = "  +3.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.227"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7112"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.424"
$ws.Range("E16").Value = "  +5.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008488"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").Value = "29.371.02"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.09"
$ws.Range("E19").Value = "  +1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").Value = "2.129.31"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.734"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1608"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.56"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.034"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.284"
$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05354"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.934"
$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7631"
$ws.Range("E35").Value = "  +1.84%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.698"
$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("D39").Value = "1.260.84"
$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.760"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.437"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.05"
$ws.Range("E42").Value = "  +3.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9044"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.09"
$ws.Range("E44").Value = "  +2.33%  "

$ws.Range("E45").Value = "  +5.66%  "

$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").Value = "2.023.96"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.799"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5196"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.474"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4338"
$ws.Range("E51").Value = "  +0.06%  "
